$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Sending cluster x Target cluster combinations (rows 2-9),
# following Dr Hou advice: re-derive ligand/receptor stats across all
# sending clusters (ECs, FAPs, M2, sCs) against target clusters (M2, sCs).
$data = @(
    ,@("ECs", "Col9a3", "Mag", "M2", 3, 1, 0.8355939999999999, 2.506782, 0.1751713572229983, 0.1751713572229983, 3, 1, 0.3133113333333333, 0.939934, 0.4010297802586483, 0.4010297802586483, 0.2618010702653333, 2.356209632388, 0.07024893089474821, 0.07024893089474821)
    ,@("ECs", "Col9a3", "Mag", "sCs", 3, 1, 0.8355939999999999, 2.506782, 0.1751713572229983, 0.1751713572229983, 2, 0.6666666666666666, 0.4679556666666667, 1.403867, 0.5989702197413518, 0.5989702197413518, 0.3910209473326666, 3.519188525994, 0.1049224263282501, 0.1049224263282501)
    ,@("FAPs", "Col9a3", "Mag", "M2", 3, 1, 2.308568333333333, 6.925705, 0.4839611679739625, 0.4839611679739625, 3, 1, 0.3133113333333333, 0.939934, 0.4010297802586483, 0.4010297802586483, 0.7233006226077777, 6.50970560347, 0.194082840846317, 0.1940828408463169)
    ,@("FAPs", "Col9a3", "Mag", "sCs", 3, 1, 2.308568333333333, 6.925705, 0.4839611679739625, 0.4839611679739625, 2, 0.6666666666666666, 0.4679556666666667, 1.403867, 0.5989702197413518, 0.5989702197413518, 1.080307633470555, 9.722768701234999, 0.2898783271276456, 0.2898783271276456)
    ,@("M2", "Col9a3", "Mag", "M2", 3, 1, 0.9553923333333335, 2.866177, 0.2002855115168937, 0.2002855115168937, 3, 1, 0.3133113333333333, 0.939934, 0.4010297802586483, 0.4010297802586483, 0.2993352458131112, 2.694017212318001, 0.08032045467261086, 0.08032045467261084)
    ,@("M2", "Col9a3", "Mag", "sCs", 3, 1, 0.9553923333333335, 2.866177, 0.2002855115168937, 0.2002855115168937, 2, 0.6666666666666666, 0.4679556666666667, 1.403867, 0.5989702197413518, 0.5989702197413518, 0.4470812562732223, 4.023731306459, 0.1199650568442829, 0.1199650568442829)
    ,@("sCs", "Col9a3", "Mag", "M2", 3, 1, 0.6705973333333333, 2.011792, 0.1405819632861455, 0.1405819632861454, 3, 1, 0.3133113333333333, 0.939934, 0.4010297802586483, 0.4010297802586483, 0.2101057446364444, 1.890951701728, 0.05637755384497228, 0.05637755384497227)
    ,@("sCs", "Col9a3", "Mag", "sCs", 3, 1, 0.6705973333333333, 2.011792, 0.1405819632861455, 0.1405819632861454, 2, 0.6666666666666666, 0.4679556666666667, 1.403867, 0.5989702197413518, 0.5989702197413518, 0.3138098221848888, 2.824288399664, 0.0842044094411732, 0.08420440944117319)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
